# Auto-generated edit script updating crypto price/volume columns (D, E) for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.019.68"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -3.13%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.195.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.83%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'571.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -2.38%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'170.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -5.58%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.615"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -5.82%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.02%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'3.195.92"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.75%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -3.04%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'6.76"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.23%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -4.40%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'3.748.61"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -2.93%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  -1.79%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'64.146.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -2.98%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'25.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -3.32%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.0000160"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.22%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.207.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -2.57%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'416.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -4.64%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -1.73%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -2.82%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -3.19%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.02%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'70.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -2.13%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +2.32%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -3.41%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -2.68%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'8.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.19%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.996"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.60%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -5.81%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'21.94"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.70%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +0.10%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'5.01"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -3.48%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'6.42"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -2.81%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -4.20%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'156.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.12%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -2.42%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'2.737.66"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -2.02%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -3.47%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'25.25"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -4.38%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -3.33%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -6.70%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'38.96"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -3.12%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'5.77"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -4.38%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.0630"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -4.68%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -4.17%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'297.99"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -6.71%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -7.80%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.0263"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.82%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -5.78%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.03%  "
$ws.Range("E51").Style = "Normal"
